# Update the install tracker: record the completed version for the
# "Roll Negative Die" feature (row 11, column C) and refresh the
# AutoFilter so rows that already have a Completed Version are
# filtered out (blank-only filter on the "Completed Version" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the version that shipped the feature.
$ws.Range("C11").Value = "1.5.3"

# (Re)apply the AutoFilter across the full data range, now that row 13
# is included, and filter column C ("Completed Version") to blanks only.
$ws.Range("A1:E13").AutoFilter()
$ws.Range("A1:E13").AutoFilter(3, @(""), 7)

# Keep the workbook-level _FilterDatabase name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$13"
    }
}

# Move the active selection, matching where the editor left off.
$ws.Range("D21").Select()
